# Adds a "Standard-Abweichung" (standard deviation) column to each of the
# seven result sheets in the Vortex Tunnel Test workbook, and restores the
# view/selection state that was active when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 - "Standard" (13 data columns B:N, average already in col O)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Standard")
$ws.Range("P1").Value = "Standard-Abweichung"
$ws.Range("P2:P9").Formula = "=STDEVA(B2:N2)"
$ws.Columns.Item(16).EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# Sheet 2 - "01 Breiter Tunnel" (7 data columns B:H, average in new col I)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("01 Breiter Tunnel")
$ws.Range("I1").Value = "Durchschnitt"
$ws.Range("J1").Value = "Standard-Abweichung"
$ws.Range("J2:J8").Formula = "=STDEVA(B2:H2)"

# ---------------------------------------------------------------------------
# Sheet 3 - "02 Enger Tunnel" (6 data columns B:G, average in new col H)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("02 Enger Tunnel")
$ws.Range("H1").Value = "Durchschnitt"
$ws.Range("I1").Value = "Standard-Abweichung"
$ws.Range("I2:I8").Formula = "=STDEVA(B2:G2)"

# ---------------------------------------------------------------------------
# Sheet 4 - "03 Schneller Tunnel" (7 data columns B:H, average in new col I)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("03 Schneller Tunnel")
$ws.Range("I1").Value = "Durchschnitt"
$ws.Range("J1").Value = "Standard-Abweichung"
$ws.Range("J2:J8").Formula = "=STDEVA(C2:H2)"

# ---------------------------------------------------------------------------
# Sheet 5 - "04 Langsamer Tunnel" (6 data columns B:G, average in new col H)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("04 Langsamer Tunnel")
$ws.Range("H1").Value = "Durchschnitt"
$ws.Range("I1").Value = "Standard-Abweichung"
$ws.Range("I2:I8").Formula = "=STDEVA(B2:G2)"

# ---------------------------------------------------------------------------
# Sheet 6 - "05 Wenige intensive Lichter" (7 data columns B:H, average in new col I)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("05 Wenige intensive Lichter")
$ws.Range("I1").Value = "Durchschnitt"
$ws.Range("J1").Value = "Standard-Abweichung"
$ws.Range("J2:J8").Formula = "=STDEVA(C2:H2)"

# ---------------------------------------------------------------------------
# Sheet 7 - "06 Viele schwache Lichter" (6 data columns B:G, average in new col H)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("06 Viele schwache Lichter")
$ws.Range("H1").Value = "Durchschnitt"
$ws.Range("I1").Value = "Standard-Abweichung"
$ws.Range("I2:I8").Formula = "=STDEVA(B2:G2)"

# ---------------------------------------------------------------------------
# Restore each sheet's selection / scroll state
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Standard")
$ws.Activate()
$ws.Range("Q4").Select()

$ws = $wb.Worksheets.Item("01 Breiter Tunnel")
$ws.Activate()
$ws.Range("K7").Select()

$ws = $wb.Worksheets.Item("02 Enger Tunnel")
$ws.Activate()
$ws.Range("K6").Select()

$ws = $wb.Worksheets.Item("03 Schneller Tunnel")
$ws.Activate()
$ws.Range("J2:J8").Select()

$ws = $wb.Worksheets.Item("04 Langsamer Tunnel")
$ws.Activate()
$ws.Range("I2:I8").Select()

$ws = $wb.Worksheets.Item("05 Wenige intensive Lichter")
$ws.Activate()
$ws.Range("J2:J8").Select()

$ws = $wb.Worksheets.Item("06 Viele schwache Lichter")
$ws.Activate()
$ws.Range("I2:I8").Select()
